# Add two new batch-sql DML test rows (batch_005, batch_006) to Sheet1,
# right after the existing batch_004 row (row 5).
#
# commit message: "exchange mysql and dingo case order, add some hash
# partition dml cases" -- the visible, reproducible part of that change in
# this worksheet is the two new rows of test-case data (rows 6 & 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6 : batch_005 -----------------------------------------------
$ws.Range("A6").Value = "batch_005"
$ws.Range("B6").Value = "y"
$ws.Range("C6").Value = "批量操作语句5执行"
$ws.Range("D6").Value = "batchsql"
$ws.Range("F6").Value = "batch05"
$ws.Range("H6").Value = "batch_sql_05"
$ws.Range("I6").Value = 'select * from $batch05'
$ws.Range("J6").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_005.csv"
$ws.Range("M6").Value = "csv_containsAll"

# ---- Row 7 : batch_006 -----------------------------------------------
$ws.Range("A7").Value = "batch_006"
$ws.Range("B7").Value = "y"
$ws.Range("C7").Value = "批量操作语句6执行"
$ws.Range("D7").Value = "batchsql"
$ws.Range("F7").Value = "batch06"
$ws.Range("H7").Value = "batch_sql_06"
$ws.Range("I7").Value = 'select * from $batch06'
$ws.Range("J7").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_006.csv"
$ws.Range("M7").Value = "csv_containsAll"

# Query_sql1 / Query_result1 columns (I & J) render right-aligned /
# "fill" in this sheet -- match the existing rows' formatting.
$ws.Range("J6").HorizontalAlignment = 5
$ws.Range("J7").HorizontalAlignment = 5

# Move the active selection to J7, matching where the author ended up
# after typing in the new rows.
$ws.Range("J7").Select()
